$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 109, shifting existing rows 109-135 down to 110-136
$ws.Rows(109).Insert()

# Populate the newly inserted row 109 with data
$ws.Range("A109").Value = 1
$ws.Range("B109").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C109").Value = 'Arica y Parinacota'
$ws.Range("D109").Value = 44627
$ws.Range("E109").Value = 15
$ws.Range("F109").Value = 'Fruta'
$ws.Range("G109").Value = 100108
$ws.Range("H109").Value = 'Tropicales y subtropicales'
$ws.Range("I109").Value = 100108002
$ws.Range("J109").Value = 'Mango'
$ws.Range("K109").Value = 'Sin especificar'
$ws.Range("L109").Value = 'Primera'
$ws.Range("M109").Value = 456
$ws.Range("N109").Value = 6500
$ws.Range("O109").Value = 7000
$ws.Range("P109").Value = 6750
$ws.Range("Q109").Value = '$/bandeja 4 kilos'
$ws.Range("R109").Value = 'Perú'
$ws.Range("S109").Value = 1688
$ws.Range("T109").Value = 4
